$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 541 (existing rows 541:549 shift down to 543:551)
$ws.Rows("541:542").Insert()

# Row 541 (new): Agrícola del Norte S.A. de Arica - Choclo - Segunda, 03-02-2022 (44595)
$ws.Cells.Item(541, 1).Value = 1
$ws.Cells.Item(541, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(541, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(541, 4).Value = 44595
$ws.Cells.Item(541, 5).Value = 15
$ws.Cells.Item(541, 6).Value = 100112024
$ws.Cells.Item(541, 7).Value = "Choclo"
$ws.Cells.Item(541, 8).Value = "Lluteño"
$ws.Cells.Item(541, 9).Value = "Segunda"
$ws.Cells.Item(541, 10).Value = 90
$ws.Cells.Item(541, 11).Value = 22000
$ws.Cells.Item(541, 12).Value = 23000
$ws.Cells.Item(541, 13).Value = 22500
$ws.Cells.Item(541, 14).Value = "$/saco 75 unidades"
$ws.Cells.Item(541, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(541, 16).Value = 300
$ws.Cells.Item(541, 17).Value = 75
$ws.Cells.Item(541, 18).Value = "Hortaliza"

# Row 542 (new): Agrícola del Norte S.A. de Arica - Choclo - Tercera, 03-02-2022 (44595)
$ws.Cells.Item(542, 1).Value = 1
$ws.Cells.Item(542, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(542, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(542, 4).Value = 44595
$ws.Cells.Item(542, 5).Value = 15
$ws.Cells.Item(542, 6).Value = 100112024
$ws.Cells.Item(542, 7).Value = "Choclo"
$ws.Cells.Item(542, 8).Value = "Lluteño"
$ws.Cells.Item(542, 9).Value = "Tercera"
$ws.Cells.Item(542, 10).Value = 70
$ws.Cells.Item(542, 11).Value = 18000
$ws.Cells.Item(542, 12).Value = 19000
$ws.Cells.Item(542, 13).Value = 18500
$ws.Cells.Item(542, 14).Value = "$/saco 100 unidades"
$ws.Cells.Item(542, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(542, 16).Value = 185
$ws.Cells.Item(542, 17).Value = 100
$ws.Cells.Item(542, 18).Value = "Hortaliza"
